$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '259.06'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '0.68%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '26.91'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-0.29%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '4.691'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '0.35%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.06042'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '2.75%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.674'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '0.72%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8594'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '0.07%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9178'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-3.24%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1395'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-0.96%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.05375'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '34.93%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07081'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-0.16%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03070'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-3.49%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09125'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.37%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001529'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.93%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0006072'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.54%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.006068'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-2.41%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.468'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-1.61%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.175'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-1.04%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.165'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-1.82%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.3127'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '2.46%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1298'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-0.19%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.134'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '8.16%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04242'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '0.13%'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-0.27%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004025'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-6.35%'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-0.08%'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-21.35%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03850'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '0.53%'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '1.18%'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-35.08%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.01513'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '32.03%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002198'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-0.08%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005158'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-5.30%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000749'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-0.08%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05450'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-22.13%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1321'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-42.78%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002098'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-0.08%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0001998'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.08%'
